# Apply the edits described by the commit "update Excel and CVS to the images/table folde"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header rename: page_path -> page_name
$ws.Range("C1").Value = "page_name"

# 2. Row 8 (US Core Implantable Device Profile -> US Core Device Profile)
$ws.Range("B8").Value = "US Core Device Profile"
$ws.Range("C8").Value = "StructureDefinition-us-core-device.html"
$ws.Range("E8").Value = "medical device"

# 3. Row 9 typo fix: clniical note -> clinical note
$ws.Range("E9").Value = "clinical note"

# 4. Append two new rows (64 and 65) with new profile data
$ws.Range("A64").Value = 61
$ws.Range("B64").Value = "US Core FamilyMemberHistory Profile"
$ws.Range("C64").Value = "StructureDefinition-us-core-familymemberhistory.html"
$ws.Range("D64").Value = "SHALL"
$ws.Range("E64").Value = "family health history"
$ws.Range("F64").Value = "FamilyMemberHistory"

$ws.Range("A65").Value = 62
$ws.Range("B65").Value = "US Core PMO ServiceRequest Profile"
$ws.Range("C65").Value = "StructureDefinition-us-core-pmo-servicerequest.html"
$ws.Range("D65").Value = "SHALL"
$ws.Range("E65").Value = "pmo order"
$ws.Range("F65").Value = "ServiceRequest"
